$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week of 11.03-15.03: fill in Friday (K/L) for "Планируемые" (row 11/12) ---
$ws.Range("K11").Value = 0.55208333333333337
$ws.Range("L11").Value = 6.5
$ws.Range("K12").Value = 0.82291666666666663

# --- Week of 11.03-15.03: fill in Friday (K/L) + sum (M) for "Фактические" (row 14/15) ---
$ws.Range("K14").Value = 0.55208333333333337
$ws.Range("L14").Value = 6.5
$ws.Range("M14").Value = 26.5
$ws.Range("K15").Value = 0.82291666666666663

# --- Week of 25.03-29.03: fill in whole week for "Планируемые" (row 27/28) ---
$ws.Range("C27").Value = 0.625
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 0.63541666666666663
$ws.Range("F27").Value = 4.5
$ws.Range("G27").Value = 0.48958333333333331
$ws.Range("H27").Value = 8
$ws.Range("I27").Value = 0.63541666666666663
$ws.Range("J27").Value = 4.5
$ws.Range("K27").Value = 0.63541666666666663
$ws.Range("L27").Value = 4.5
$ws.Range("M27").Value = 26.5

$ws.Range("C28").Value = 0.83333333333333337
$ws.Range("E28").Value = 0.82291666666666663
$ws.Range("G28").Value = 0.82291666666666663
$ws.Range("I28").Value = 0.82291666666666663
$ws.Range("K28").Value = 0.82291666666666663

# --- Week of 25.03-29.03: fill in part of the week for "Фактические" (row 30/31), today = Thursday ---
$ws.Range("C30").Value = 0.625
$ws.Range("D30").Value = 5
$ws.Range("E30").Value = 0.63541666666666663
$ws.Range("F30").Value = 4.5
$ws.Range("G30").Value = 0.47916666666666669
$ws.Range("H30").Value = 8

$ws.Range("C31").Value = 0.83333333333333337
$ws.Range("E31").Value = 0.82291666666666663
$ws.Range("G31").Value = 0.8125

# --- Update "today" selection to the last filled cell (H30:H31, today's hours) ---
$ws.Range("H30:H31").Select()
